$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("cost_A", 37.689441000000002, 129.03623099999999),
    @("cost_B", 37.640802000000001, 129.044299),
    @("cost_C", 37.599268000000002, 129.082065),
    @("cost_D", 37.556690000000003, 129.12120400000001),
    @("cost_E", 37.511091999999998, 129.12944400000001),
    @("cost_F", 37.402901, 129.21115499999999),
    @("cost_G", 37.342744000000003, 129.26574299999999),
    @("cost_H", 37.160060000000001, 129.355694),
    @("cost_I", 36.980088000000002, 129.40719200000001)
)

$row = 10
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

$ws.Range("E13").Select()
